# Update countries & provincias Spain
# - Update "Datos actualizados ..." timestamp in A1
# - España overtakes Estados Unidos (new data for España, Estados Unidos keeps
#   its previous figures but drops one rank)
# - Noruega stats refreshed in place
# - Islandia jumps ahead of Chile with fresh data; Chile/Polonia/Ecuador/
#   Grecia/Catar each shift down one rank (figures unchanged per country)
# - Hungria stats refreshed in place

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 12:16"

# Row 6: España (new totals, moves above Estados Unidos)
$ws.Range("A6").Value = "España"
$ws.Range("B6").Value = 28572
$ws.Range("C6").Value = 3076
$ws.Range("D6").Value = 2125
$ws.Range("E6").Value = 25066
$ws.Range("F6").Value = 1612
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1381

# Row 7: Estados Unidos (keeps its previous totals, now one rank lower)
$ws.Range("A7").Value = "Estados Unidos"
$ws.Range("B7").Value = 26900
$ws.Range("C7").Value = 2693
$ws.Range("D7").Value = 178
$ws.Range("E7").Value = 26374
$ws.Range("F7").Value = 708
$ws.Range("G7").Value = 46
$ws.Range("H7").Value = 348

# Row 17: Noruega (figures refreshed in place)
$ws.Range("A17").Value = "Noruega"
$ws.Range("B17").Value = 2219
$ws.Range("C17").Value = 55
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 2206
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 7

# Row 35: Islandia (new totals, moves above Chile)
$ws.Range("A35").Value = "Islandia"
$ws.Range("B35").Value = 568
$ws.Range("C35").Value = 95
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 562
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 1

# Row 36: Chile (keeps its previous totals, now one rank lower)
$ws.Range("A36").Value = "Chile"
$ws.Range("B36").Value = 537
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 8
$ws.Range("E36").Value = 528
$ws.Range("F36").Value = 7
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 1

# Row 37: Polonia (keeps its previous totals, now one rank lower)
$ws.Range("A37").Value = "Polonia"
$ws.Range("B37").Value = 536
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 13
$ws.Range("E37").Value = 518
$ws.Range("F37").Value = 3
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 5

# Row 38: Ecuador (keeps its previous totals, now one rank lower)
$ws.Range("A38").Value = "Ecuador"
$ws.Range("B38").Value = 532
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 3
$ws.Range("E38").Value = 522
$ws.Range("F38").Value = 2
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 7

# Row 39: Grecia (keeps its previous totals, now one rank lower)
$ws.Range("A39").Value = "Grecia"
$ws.Range("B39").Value = 530
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 19
$ws.Range("E39").Value = 498
$ws.Range("F39").Value = 18
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 13

# Row 40: Catar (keeps its previous totals, now one rank lower)
$ws.Range("A40").Value = "Catar"
$ws.Range("B40").Value = 481
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 27
$ws.Range("E40").Value = 454
$ws.Range("F40").Value = 6
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0

# Row 72: Hungria (figures refreshed in place)
$ws.Range("A72").Value = "Hungria"
$ws.Range("B72").Value = 131
$ws.Range("C72").Value = 28
$ws.Range("D72").Value = 16
$ws.Range("E72").Value = 109
$ws.Range("F72").Value = 6
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 6
